# Add a "Save" column (H) to the s_vals sheet, mirroring the header
# formatting already used by the other header cells (e.g. G1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: text "Save", formatted like the rest of row 1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells H2:H3 = 1 (unformatted, like the other numeric cells).
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
